# Mediciones Derivador Compensado - "Almost finished, missing the integrator"
#
# This script:
#  1. Fixes four measured "C" values on the "Respuesta en frecuencia" sheet
#     (rows that were mis-measured/updated), then sorts the A3:F23 data
#     table by frequency (column B) ascending - matching the interleaved
#     extra measurements (15000/20000/150000/180000/300000 Hz) being
#     merged back into the sorted sweep.
#  2. Tweaks a handful of phase/gain inputs on the "Impedancia de entrada"
#     sheet (columns F/G) which ripple through the H/I/J formula columns.
#  3. Restores the selection/active-sheet state left by the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: "Impedancia de entrada"
# ---------------------------------------------------------------------
$wsImp = $wb.Worksheets.Item("Impedancia de entrada")

$wsImp.Range("G4").Value = 50
$wsImp.Range("G9").Value = 30
$wsImp.Range("F12").Value = 0.36
$wsImp.Range("G12").Value = 80
$wsImp.Range("F15").Value = 0.2

$wsImp.Range("G13").Select()

# ---------------------------------------------------------------------
# Sheet: "Respuesta en frecuencia"
# ---------------------------------------------------------------------
$wsResp = $wb.Worksheets.Item("Respuesta en frecuencia")

# Corrected measurements (still in their original, unsorted rows)
$wsResp.Range("C15").Value = 14.9
$wsResp.Range("C21").Value = 7.58
$wsResp.Range("C22").Value = 10.26
$wsResp.Range("C16").Value = 9.11

# Re-sort the measurement table by frequency (column B) ascending, now
# that the extra points (15000/20000/150000/180000/300000 Hz) should be
# interleaved in order with the rest of the sweep.
$sortRange = $wsResp.Range("A3:F23")
$keyRange = $wsResp.Range("B3:B23")
$sortRange.Sort($keyRange, 1)

$wsResp.Range("A13").Select()
$wsResp.Activate()
